$wb = $excel.ActiveWorkbook

# --- Users sheet ----------------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Julie Carthane"
$wsUsers.Range("A2").Style = "Normal"
$wsUsers.Range("B2").Value = "Aja Mount"
$wsUsers.Range("A2").Select() | Out-Null

# --- AppName sheet ----------------------------------------------------------
# (no content/selection change on this sheet)

# --- ModuleName sheet -------------------------------------------------------
$wsModuleName = $wb.Worksheets.Item("ModuleName")
$wsModuleName.Range("D20").Select() | Out-Null

# --- GiftLog sheet ----------------------------------------------------------
$wsGiftLog = $wb.Worksheets.Item("GiftLog")
$wsGiftLog.Range("B2").Value = "Julie Carthane"
$wsGiftLog.Range("B2").Style = "Normal"

# GiftLog becomes the active/selected tab, with D9 as the active cell
$wsGiftLog.Activate() | Out-Null
$wsGiftLog.Range("D9").Select() | Out-Null
